# Draup-Automation FW update-SignIN,Universe-Done Required Changes-Renuka-10-11-17.
#
# 1) Insert a new "Sheet1" worksheet between "Elements" and "LoginData",
#    carrying the "UN=new UniversePageObjects(driver, test)" snippet in E1,
#    and make it the active tab.
# 2) On "Elements": rename the "Universe" locator row to "Universe,Opportunity",
#    and add a batch of new locator rows (table headers, vertical filters,
#    signal columns, subverticals) below the existing data, formatting the
#    new "vertical" rows (automotive/consumersoftware/enterprisesoftware/
#    subverticals) with a bold-free black font + wrapped/centered alignment.
# 3) Resize columns C/D on "Elements" to fit the new, longer xpath text.

$wb = $excel.ActiveWorkbook
$elements = $wb.Worksheets.Item("Elements")

# --- 1. Insert the new "Sheet1" worksheet right after "Elements" ---------
$newSheet = $wb.Worksheets.Add($null, $elements)
$newSheet.Name = "Sheet1"
$newSheet.Range("E1").Value = "  UN=new UniversePageObjects(driver, test)"

# --- 2. Extend "Elements" with the new locator rows -----------------------
$ws = $elements

# 2a. Baseline formatting (font/fill/border already used by row 7) across
#     the full new block of rows, same as the existing data rows below the
#     header.
$ws.Range("A7:C7").Copy()
$ws.Range("A8:C21").PasteSpecial(-4122)

# 2b. Build the new "vertical label" style (black font, wrapped, vertically
#     centered) once on a scratch cell far away from the used range, then
#     stamp it onto the real target cells and discard the scratch cell.
$ws.Range("A7").Copy()
$ws.Range("Z100").PasteSpecial(-4122)
$scratch = $ws.Range("Z100")
$scratch.Font.Color = 0
$scratch.WrapText = $true
$scratch.VerticalAlignment = -4108
$scratch.Copy()
$ws.Range("A17:C20").PasteSpecial(-4122)
$ws.Range("Z100").Delete()

# 2c. Column B of those rows keeps the plain existing style (it's still an
#     "xpath" label cell like all the others).
$ws.Range("B7:B7").Copy()
$ws.Range("B17:B20").PasteSpecial(-4122)

# 2d. Cell values.
$ws.Cells.Item(6, 1).Value = "Universe,Opportunity"

$ws.Cells.Item(8, 1).Value = "TableAccoutHeader"
$ws.Cells.Item(8, 2).Value = "xpath"
$ws.Cells.Item(8, 3).Value = "//div[@class='text-left ag-header-cell-label'][@title='Account']"

$ws.Cells.Item(9, 1).Value = "TableOpportunityIndex"
$ws.Cells.Item(9, 2).Value = "xpath"
$ws.Cells.Item(9, 3).Value = "//div[@class='text-center ag-header-cell-label'][@title='Opportunity Index']"

$ws.Cells.Item(10, 1).Value = "BaseOpportunityIndex"
$ws.Cells.Item(10, 2).Value = "xpath"
$ws.Cells.Item(10, 3).Value = "//div[@class='text-center ag-header-cell-label'][@title='Base Opportunity Index']"

$ws.Cells.Item(11, 1).Value = "IncrementalOpportunity"
$ws.Cells.Item(11, 2).Value = "xpath"
$ws.Cells.Item(11, 3).Value = "//div[@class='text-center ag-header-cell-label'][@title='Incremental Opportunity Index']"

$ws.Cells.Item(12, 1).Value = "HiringSignal"
$ws.Cells.Item(12, 2).Value = "xpath"
$ws.Cells.Item(12, 3).Value = "//div[@class='text-center ag-header-cell-label'][@title='Hiring Signal']"

$ws.Cells.Item(13, 1).Value = "InvestementSignals"
$ws.Cells.Item(13, 2).Value = "xpath"
$ws.Cells.Item(13, 3).Value = "//div[@class='text-center ag-header-cell-label'][@title='Investments Signal']"

$ws.Cells.Item(14, 1).Value = "M&Asignal"
$ws.Cells.Item(14, 2).Value = "xpath"
$ws.Cells.Item(14, 3).Value = "//div[@class='text-center ag-header-cell-label'][@title='M&A Signal']"

$ws.Cells.Item(15, 1).Value = "LayoffSignal"
$ws.Cells.Item(15, 2).Value = "xpath"
$ws.Cells.Item(15, 3).Value = "//div[@class='text-center ag-header-cell-label'][@title='Layoff/Attrition Signal']"

$ws.Cells.Item(16, 1).Value = "ExecutiveMovement"
$ws.Cells.Item(16, 3).Value = "//div[@class='text-center ag-header-cell-label'][@title='Executive Movement Signal']"

$ws.Cells.Item(17, 1).Value = "automotive"
$ws.Cells.Item(17, 2).Value = "xpath"
$ws.Cells.Item(17, 3).Value = "//label[text()='Automotive']"

$ws.Cells.Item(18, 1).Value = "consumersoftware"
$ws.Cells.Item(18, 2).Value = "xpath"
$ws.Cells.Item(18, 3).Value = "//label[text()='Consumer Software']"

$ws.Cells.Item(19, 1).Value = "enterprisesoftware"
$ws.Cells.Item(19, 2).Value = "xpath"
$ws.Cells.Item(19, 3).Value = "//label[text()='Enterprise Software']"

$ws.Cells.Item(20, 1).Value = "subverticals"
$ws.Cells.Item(20, 2).Value = "xpath"
$ws.Cells.Item(20, 3).Value = "//h6[text()='Subverticals']"

# --- 3. Widen columns C/D to fit the new (longer) xpath text --------------
$ws.Columns.Item(3).ColumnWidth = 81
$ws.Columns.Item(4).ColumnWidth = 22

# --- 4. Selection / active-tab state ---------------------------------------
$newSheet.Range("E1").Select()
$ws.Range("C6").Select()
$newSheet.Activate()
